# Update the "Marking" row (row 11) and "Total" row (row 12) of the
# concise_ms marksheet pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right 4 -> 5, Wrong -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 ("Total"): Right 56 -> 70, Wrong 0 -> -0, Max text "56/112" -> "70.0/140"
$ws.Range("B12").Value = 70
$ws.Range("C12").Value = -0
$ws.Range("E12").Value = "70.0/140"
